# ADD results from server
# Update computed result values (row 2) on each year sheet with the
# latest figures returned by the server.

$wb = $excel.ActiveWorkbook

function Set-Results {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Values.Keys) {
        $ws.Range([string]$cellRef).Value = $Values[$cellRef]
    }
}

Set-Results "2025" @{
    "A2" = 3906.399109145206
    "C2" = 48353.76274462014
    "F2" = 9433.134471502228
    "H2" = 2534.277928792104
    "N2" = 2368.085410105587
    "O2" = 1996.112248849888
}

Set-Results "2030" @{
    "A2" = 6991.052031681918
    "C2" = 197913.7502057619
    "F2" = 16452.51445364119
    "H2" = 8194.52068131253
    "N2" = 7546.08779474431
    "O2" = 6258.805786094594
}

Set-Results "2035" @{
    "A2" = 31236.29455387744
    "C2" = 292247.2772138842
    "F2" = 16595.10705160327
    "H2" = 12131.91920790125
    "N2" = 12893.74290450278
    "O2" = 9265.399564033123
}

Set-Results "2040" @{
    "A2" = 31236.29455387744
    "C2" = 292247.2772138842
    "F2" = 16595.10705160327
    "H2" = 12131.91920790125
    "N2" = 14051.59265694596
    "O2" = 9265.399564033123
}

Set-Results "2045" @{
    "A2" = 38906.8534480406
    "B2" = 193.0947398408091
    "C2" = 292247.2772138842
    "F2" = 16595.10705160327
    "H2" = 12131.91920790125
    "N2" = 16889.45207698157
    "O2" = 10098.21437025084
}

Set-Results "2050" @{
    "A2" = 38906.8534480406
    "B2" = 193.0947398408091
    "C2" = 292247.2772138842
    "F2" = 16595.10705160327
    "H2" = 12131.91920790125
    "N2" = 16889.45207698157
    "O2" = 10098.21437025084
}
